# Update countries & provincias Spain
# Applies the "22:30" data refresh on top of the "21:13" snapshot:
#  - swaps the Siria / Sri Lanka rows (order + values) so that row 135
#    becomes Siria's updated figures and row 136 becomes Sri Lanka's
#    (unchanged) figures
#  - refreshes numeric statistics for a handful of countries
#  - updates the "Datos actualizados..." timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Octubre de 2020 a las 22:30"

# --- Helper: write a full data row (columns A-H), positional args -------
function Set-Row {
    param($Row, $Country, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes)
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# --- Row 4: Estados Unidos ----------------------------------------------
Set-Row 4 "Estados Unidos" 7931821 37343 5080869 2631884 0 421 219068

# --- Row 5: India ---------------------------------------------------------
Set-Row 5 "India" 7051543 74535 6074863 868309 0 921 108371

# --- Row 13: Francia (only Recuperados / Muertes hoy / Muertes change) ---
$ws.Cells.Item(13, 5).Value = 585408
$ws.Cells.Item(13, 7).Value = 54
$ws.Cells.Item(13, 8).Value = 32637

# --- Row 14: Sudafrica -----------------------------------------------------
Set-Row 14 "Sudafrica" 690896 2544 622153 51070 0 126 17673

# --- Row 15: Reino Unido (Casos activos/Recuperados stay 0) ---------------
$ws.Cells.Item(15, 2).Value = 590844
$ws.Cells.Item(15, 3).Value = 15166
$ws.Cells.Item(15, 7).Value = 81
$ws.Cells.Item(15, 8).Value = 42760

# --- Row 25: Alemania ------------------------------------------------------
$ws.Cells.Item(25, 2).Value = 323438
$ws.Cells.Item(25, 3).Value = 2960
$ws.Cells.Item(25, 5).Value = 40247
$ws.Cells.Item(25, 7).Value = 4
$ws.Cells.Item(25, 8).Value = 9691

# --- Row 27: Israel ----------------------------------------------------
Set-Row 27 "Israel" 289875 2017 225725 62209 0 55 1941

# --- Row 50: Costa Rica ----------------------------------------------------
Set-Row 50 "Costa Rica" 87439 1386 52669 33694 0 21 1076

# --- Row 53: Etiopia ---------------------------------------------------
Set-Row 53 "Etiopia" 83429 767 37683 44469 0 6 1277

# --- Row 90: Costa de Marfil (F/G/H unchanged) -----------------------------
$ws.Cells.Item(90, 2).Value = 20128
$ws.Cells.Item(90, 3).Value = 92
$ws.Cells.Item(90, 4).Value = 19752
$ws.Cells.Item(90, 5).Value = 256

# --- Row 109: Mozambique -----------------------------------------------
Set-Row 109 "Mozambique" 9844 102 7203 2571 0 1 70

# --- Row 120: Angola -----------------------------------------------------
Set-Row 120 "Angola" 6246 215 2716 3312 0 6 218

# --- Row 124: Suazilandia (F/G/H unchanged) --------------------------------
$ws.Cells.Item(124, 2).Value = 5660
$ws.Cells.Item(124, 3).Value = 16
$ws.Cells.Item(124, 4).Value = 5282
$ws.Cells.Item(124, 5).Value = 265

# --- Row 133: Ruanda (F/G/H unchanged) -------------------------------------
$ws.Cells.Item(133, 2).Value = 4892
$ws.Cells.Item(133, 3).Value = 2
$ws.Cells.Item(133, 4).Value = 3566
$ws.Cells.Item(133, 5).Value = 1296

# --- Rows 135/136: Siria & Sri Lanka swap places ---------------------------
# Row 135 now shows Siria's refreshed figures, row 136 shows Sri Lanka's
# (unchanged since the last snapshot) figures.
Set-Row 135 "Siria" 4673 57 1271 3181 0 3 221
Set-Row 136 "Sri Lanka" 4628 105 3306 1309 0 0 13

# --- Row 144: Gambia (F/G/H unchanged) -------------------------------------
$ws.Cells.Item(144, 2).Value = 3628
$ws.Cells.Item(144, 3).Value = 7
$ws.Cells.Item(144, 4).Value = 2540
$ws.Cells.Item(144, 5).Value = 971

# --- Row 165: Liberia (D unchanged) ----------------------------------------
$ws.Cells.Item(165, 2).Value = 1363
$ws.Cells.Item(165, 3).Value = 3
$ws.Cells.Item(165, 5).Value = 36

# --- Row 166: Republica del Chad -----------------------------------------
Set-Row 166 "Republica del Chad" 1291 17 1103 96 0 2 92
